$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F315").Value = 57063
$ws.Range("G315").Value = 2635
$ws.Range("F450").Value = 91773
$ws.Range("F483").Value = 66737
$ws.Range("F543").Value = 4818
$ws.Range("F547").Value = 14247
$ws.Range("F555").Value = 21869
$ws.Range("F572").Value = 33867
$ws.Range("F576").Value = 29588
$ws.Range("F621").Value = 56741
$ws.Range("G621").Value = 4171
$ws.Range("F624").Value = 51753
$ws.Range("F625").Value = 43919
$ws.Range("F626").Value = 20417
$ws.Range("G626").Value = 2110
$ws.Range("F627").Value = 34304
$ws.Range("G627").Value = 2765
$ws.Range("F628").Value = 65131
$ws.Range("G628").Value = 4221
$ws.Range("F629").Value = 46514
$ws.Range("F630").Value = 47310
$ws.Range("G630").Value = 3008
$ws.Range("F631").Value = 42428
$ws.Range("F632").Value = 44556
$ws.Range("F633").Value = 24338
$ws.Range("G633").Value = 1963
$ws.Range("F634").Value = 47575
$ws.Range("G634").Value = 2213
$ws.Range("F635").Value = 84146
$ws.Range("G635").Value = 3749
$ws.Range("F636").Value = 50370
$ws.Range("G636").Value = 2365
$ws.Range("F637").Value = 43931
$ws.Range("G637").Value = 2122
$ws.Range("F638").Value = 38038
$ws.Range("G638").Value = 1977
$ws.Range("F639").Value = 40867
$ws.Range("G639").Value = 1978
$ws.Range("F640").Value = 19993
$ws.Range("G640").Value = 1241
$ws.Range("F641").Value = 34733
$ws.Range("G641").Value = 1398
$ws.Range("F642").Value = 67946
$ws.Range("G642").Value = 2409
$ws.Range("F643").Value = 43665
$ws.Range("G643").Value = 1663
$ws.Range("F644").Value = 37084
$ws.Range("G644").Value = 1507
$ws.Range("F645").Value = 35890
$ws.Range("G645").Value = 1316
$ws.Range("F646").Value = 36202
$ws.Range("F647").Value = 16354
$ws.Range("G647").Value = 916
$ws.Range("F648").Value = 30684
$ws.Range("G648").Value = 1073
$ws.Range("F649").Value = 63097
$ws.Range("G649").Value = 1820
$ws.Range("F650").Value = 38242
$ws.Range("F651").Value = 37278
$ws.Range("F652").Value = 35266
$ws.Range("F653").Value = 34295
$ws.Range("G653").Value = 1016
$ws.Range("F655").Value = 25499
$ws.Range("G655").Value = 813
$ws.Range("F656").Value = 53067
$ws.Range("F657").Value = 34226
$ws.Range("G657").Value = 880
$ws.Range("F658").Value = 27352
$ws.Range("F659").Value = 26492
$ws.Range("F662").Value = 12926
$ws.Range("G662").Value = 544
$ws.Range("F663").Value = 37773
$ws.Range("F664").Value = 26588
$ws.Range("F665").Value = 28372
$ws.Range("F666").Value = 24120
$ws.Range("G666").Value = 783
$ws.Range("H897").Value = 526
$ws.Range("H898").Value = 500
$ws.Range("H899").Value = 522
$ws.Range("H900").Value = 556
$ws.Range("H901").Value = 555
$ws.Range("H902").Value = 517
$ws.Range("H903").Value = 474
$ws.Range("H904").Value = 480
$ws.Range("H905").Value = 443
$ws.Range("H906").Value = 445
$ws.Range("H907").Value = 487
$ws.Range("H908").Value = 504
$ws.Range("H909").Value = 457
$ws.Range("H910").Value = 412
$ws.Range("H911").Value = 421
$ws.Range("H912").Value = 383
$ws.Range("H913").Value = 402
$ws.Range("H914").Value = 420
$ws.Range("H915").Value = 404
$ws.Range("H916").Value = 395
$ws.Range("H919").Value = 375
$ws.Range("H921").Value = 397
$ws.Range("H922").Value = 371
$ws.Range("H923").Value = 369
$ws.Range("H924").Value = 363
$ws.Range("H925").Value = 381
$ws.Range("H926").Value = 375
$ws.Range("H927").Value = 408
$ws.Range("H928").Value = 445
$ws.Range("H929").Value = 450
$ws.Range("H930").Value = 435
$ws.Range("H931").Value = 438
$ws.Range("H932").Value = 462
$ws.Range("F933").Value = 3064
$ws.Range("G933").Value = 369
$ws.Range("H933").Value = 462
$ws.Range("H934").Value = 514
$ws.Range("F935").Value = 1117
$ws.Range("G935").Value = 145
$ws.Range("H935").Value = 550
$ws.Range("F936").Value = 5701
$ws.Range("G936").Value = 665
$ws.Range("H936").Value = 573
$ws.Range("F937").Value = 3650
$ws.Range("G937").Value = 470
$ws.Range("H937").Value = 560
$ws.Range("F938").Value = 3648
$ws.Range("G938").Value = 426
$ws.Range("H938").Value = 574
$ws.Range("F939").Value = 4613
$ws.Range("G939").Value = 432
$ws.Range("H939").Value = 580
$ws.Range("F940").Value = 2778
$ws.Range("G940").Value = 388
$ws.Range("H940").Value = 578
$ws.Range("H941").Value = 578
$ws.Range("H942").Value = 578
